$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "dia" (row 3) and "kamu" (row 4) present with timestamps
$ws.Range("B3").Value = "HADIR"
$ws.Range("C3").Value = 45083.294571759259

$ws.Range("B4").Value = "HADIR"
$ws.Range("C4").Value = 45083.298587962963

# Remove the now-unneeded styled but empty cell for "z" (row 6)
$ws.Range("C6").Clear()

# Delete the row for "anumuhafidz" (row 7)
$ws.Rows.Item(7).Delete()

# Update the selection to match the recorded view state
$ws.Range("C9").Select()
